$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.688.76'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +4.13%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.418.65'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.54%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.03'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +4.33%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.17'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +6.15%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.513'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.21%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.529'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +11.27%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.32'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.84%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0798'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.68%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.72'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.54%  '

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.12%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.91'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.12%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.799.10'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.70%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.423.40'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.96%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.830'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.48%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.511.85'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.64%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.31'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.21%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.36'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.67%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0914'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.35%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.69'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.16%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '242.48'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.09%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.27'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.24%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.18%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.07%  '

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.28%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.19'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -6.88%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.47'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.70%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.62'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.75%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '48.63'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.60%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +17.62%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.50'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +11.55%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.15'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.98%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0771'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +6.52%  '

$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.35%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.88'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.44%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.46'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.11%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.84'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.09%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '122.99'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.25%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.78%  '

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.92%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.88'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.08%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0290'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.21%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.938.72'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.62%  '

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.41%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.93'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +7.87%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.23'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.05%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.73'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +15.09%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.59'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.90%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '53.80'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.26%  '
